# Fixed update to excel issue
# - Rename "Requested quantity" headers on existing sheets
# - Add a new "PO Forecast" sheet with forecast data

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$ws2 = $wb.Worksheets.Item(2)   # "Monthly Trend"

# Rename the "Requested quantity" column headers
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "PO Forecast"

# Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Reuse the same header formatting (bold, bordered, centered) as the other sheets
$ws1.Range("A1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122) # xlPasteFormats

# Forecast data rows
$ws3.Range("A2").Value = 45620.99999999999
$ws3.Range("B2").Value = 6
$ws3.Range("C2").Value = 6.000004066545454
$ws3.Range("D2").Value = 6.000004066628178

$ws3.Range("A3").Value = 45627.99999999999
$ws3.Range("B3").Value = 10
$ws3.Range("C3").Value = 10.00000405500016
$ws3.Range("D3").Value = 10.00000405507793

$ws3.Range("A4").Value = 45634.99999999999
$ws3.Range("B4").Value = 14
$ws3.Range("C4").Value = 14.00000395667412
$ws3.Range("D4").Value = 14.00000411703857

$ws3.Range("A5").Value = 45641.99999999999
$ws3.Range("B5").Value = 18
$ws3.Range("C5").Value = 18.00000375804334
$ws3.Range("D5").Value = 18.00000427243518

$ws3.Range("A6").Value = 45648.99999999999
$ws3.Range("B6").Value = 22
$ws3.Range("C6").Value = 22.00000348895639
$ws3.Range("D6").Value = 22.00000448120489

$ws3.Range("A7").Value = 45655.99999999999
$ws3.Range("B7").Value = 26
$ws3.Range("C7").Value = 26.00000318687867
$ws3.Range("D7").Value = 26.00000475864504

$ws3.Range("A8").Value = 45662.99999999999
$ws3.Range("B8").Value = 30
$ws3.Range("C8").Value = 30.00000286514124
$ws3.Range("D8").Value = 30.00000507659395

$ws3.Range("A9").Value = 45669.99999999999
$ws3.Range("B9").Value = 34
$ws3.Range("C9").Value = 34.00000246929502
$ws3.Range("D9").Value = 34.00000544024565

$ws3.Range("A10").Value = 45676.99999999999
$ws3.Range("B10").Value = 38
$ws3.Range("C10").Value = 38.00000206927972
$ws3.Range("D10").Value = 38.00000576646445

$ws3.Range("A11").Value = 45683.99999999999
$ws3.Range("B11").Value = 42
$ws3.Range("C11").Value = 42.00000163747394
$ws3.Range("D11").Value = 42.00000621959237

# Apply the same date/time number formatting used in column A of the other sheets
$ws1.Range("A2").Copy()
$ws3.Range("A2:A11").PasteSpecial(-4122) # xlPasteFormats

[void]$ws3.Range("A1").Select()
